$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warrant Issued Report")

# Insert a new row above row 50 (old row 50 and everything below shifts down by one)
$ws.Range("A50:F50").Insert()

# Copy the formatting of the row above (row 49) into the newly inserted row 50
$ws.Range("A49:F49").Copy()
$ws.Range("A50:F50").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row's content
$ws.Cells.Item(50, 3).Value = "Vehicle RegistrationState"
$ws.Cells.Item(50, 5).Value = "MI"
$ws.Cells.Item(50, 6).Value = "wir-doc:WarrantIssuedReport/j:ConveyanceRegistration[@structures:id=/wir-doc:WarrantIssuedReport/j:ConveyanceRegistrationAssociation/j:ItemRegistration/@structures:ref]/j:JurisdictionNCICLISCode"

Write-Host ("Dimension after: " + $ws.UsedRange.Address())
